$wb = $excel.ActiveWorkbook

# The three existing sheets, in tab order.
$wsFaculty = $wb.Worksheets.Item(1)   # "FacultyMembers" -> becomes "UniversityMembers"
$wsStudents = $wb.Worksheets.Item(2)  # "Students"
$wsStaff = $wb.Worksheets.Item(3)     # "NonFacultyStaff"

# Add the new "Type" header/column to the surviving sheet.
$wsFaculty.Range("C1").Value = "Type"
$wsFaculty.Range("C2").Value = "FacultyMember"
$wsFaculty.Range("C3").Value = "FacultyMember"

# Bring the Students rows over (rows 4-7), preserving hyperlinked emails.
$wsFaculty.Range("A4").Value = "Sheida Khodabakhshian"
$wsFaculty.Range("B4").Value = "sheida@my.yorku.ca"
$wsFaculty.Hyperlinks.Add($wsFaculty.Range("B4"), "mailto:sheida@my.yorku.ca")
$wsFaculty.Range("B4").Style = "Hyperlink"
$wsFaculty.Range("C4").Value = "Student"

$wsFaculty.Range("A5").Value = "Saaram Mustafa"
$wsFaculty.Range("B5").Value = "saaram@my.yorku.ca"
$wsFaculty.Hyperlinks.Add($wsFaculty.Range("B5"), "mailto:saaram@my.yorku.ca")
$wsFaculty.Range("B5").Style = "Hyperlink"
$wsFaculty.Range("C5").Value = "Student"

$wsFaculty.Range("A6").Value = "Abdullah Syed Ali "
$wsFaculty.Range("B6").Value = "abdullah@my.yorku.ca"
$wsFaculty.Hyperlinks.Add($wsFaculty.Range("B6"), "mailto:abdullah@my.yorku.ca")
$wsFaculty.Range("B6").Style = "Hyperlink"
$wsFaculty.Range("C6").Value = "Student"

$wsFaculty.Range("A7").Value = "Joanne Yu"
$wsFaculty.Range("B7").Value = "joanne@my.yorku.ca"
$wsFaculty.Hyperlinks.Add($wsFaculty.Range("B7"), "mailto:joanne@my.yorku.ca")
$wsFaculty.Range("B7").Style = "Hyperlink"
$wsFaculty.Range("C7").Value = "Student"

# Bring the NonFacultyStaff rows over (rows 8-9).
$wsFaculty.Range("A8").Value = "Harry Kim"
$wsFaculty.Range("B8").Value = "harry@my.yorku.ca"
$wsFaculty.Hyperlinks.Add($wsFaculty.Range("B8"), "mailto:harry@my.yorku.ca")
$wsFaculty.Range("B8").Style = "Hyperlink"
$wsFaculty.Range("C8").Value = "NonFacultyStaff"

$wsFaculty.Range("A9").Value = "Natalie Nate"
$wsFaculty.Range("B9").Value = "natalie@my.yorku.ca"
$wsFaculty.Hyperlinks.Add($wsFaculty.Range("B9"), "mailto:natalie@my.yorku.ca")
$wsFaculty.Range("B9").Style = "Hyperlink"
$wsFaculty.Range("C9").Value = "NonFacultyStaff"

# Drop the now-redundant sheets (delete from the end so indices stay valid).
$wsStaff.Delete()
$wsStudents.Delete()

# Rename the remaining sheet and select it / the next free row.
$wsFaculty.Name = "UniversityMembers"
$wsFaculty.Activate()
$wsFaculty.Range("E10").Select()
